$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-17.31100322206613, 1.336374393439415, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.562033789318932, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, 2.404350525695407, 1.569527944419806, -17.31100322206613, 3.625132176139151, -17.31100322206613, 1.637378269386086, -17.31100322206613, 2.067310686968454, -17.31100322206613),
    @(-17.31100322206613, 1.590302689376191, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.845890875467573, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(2.583064027666959, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, 4.321919666040024, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(3.808156388172628, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 1.073987050901758, -17.31100322206613, 1.995523568232192),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.911698538324045, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.100667058039246),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 1.831355666476085, 1.871965724308396),
    @(-17.31100322206613, -17.31100322206613, 0.7650868793363435, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.032904649276443),
    @(-17.31100322206613, -17.31100322206613, 1.13956771284284, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.013302601900895, -17.31100322206613),
    @(-17.31100322206613, 2.590625301038369, 2.940808544996305, -17.31100322206613, -17.31100322206613, -17.31100322206613, 1.440661780049489, 2.295497780574576, 2.241435108030472, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, -17.31100322206613, 1.87815775526138, 1.383521306615711, 1.800919474867873, -17.31100322206613),
    @(-17.31100322206613, -17.31100322206613, 1.519279860492104, -17.31100322206613, -17.31100322206613, -17.31100322206613, 1.507719368702015, 1.377250502534217, -17.31100322206613, -17.31100322206613),
    @(-17.31100322206613, 0.5637200325534802, 1.369785512998557, -17.31100322206613, 2.937546336175596, -17.31100322206613, 1.789569294715534, 0.9055579305474034, -17.31100322206613, 1.989369449044808),
    @(-17.31100322206613, 0.7381061396315566, -17.31100322206613, -17.31100322206613, -17.31100322206613, 2.402002282246908, 2.069595236813325, -17.31100322206613, -17.31100322206613, -17.31100322206613)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
